$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28:123 down to 29:124.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new record (same shape as the row that used
# to be there, with updated D/J/K/L/M/P values).
$ws.Cells.Item(28, 1).Value = 10
$ws.Cells.Item(28, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(28, 3).Value = "La Araucanía"
$ws.Cells.Item(28, 4).Value = 44600
$ws.Cells.Item(28, 5).Value = 9
$ws.Cells.Item(28, 6).Value = 100112012
$ws.Cells.Item(28, 7).Value = "Espinaca"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 95
$ws.Cells.Item(28, 11).Value = 12000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 12000
$ws.Cells.Item(28, 14).Value = "`$/docena de atados"
$ws.Cells.Item(28, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(28, 16).Value = 4000
$ws.Cells.Item(28, 17).Value = 3
$ws.Cells.Item(28, 18).Value = "Hortaliza"
